# Sample Project (DESIGN/rules) - Main.xlsx
# B11 on the "Rules" sheet currently holds the shared string "R40".
# It is changed to the text "1" (kept as text, not a number), which
# Excel appends as a brand-new shared-string table entry.
#
# A plain `.Value = "1"` assignment would let Excel auto-detect the
# numeric-looking literal and store it as a Number, so instead we
# compute it with TEXT() and then freeze the formula result down to a
# static value via Copy + PasteSpecial (values only). That keeps the
# cell's underlying type as text while leaving its style/number format
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$cell = $ws.Range("B11")
$cell.Formula = '=TEXT(1,"0")'
$cell.Copy() | Out-Null
$cell.PasteSpecial(-4163) | Out-Null
